# Devs.xlsx — "steps towards working mvp"
# Update the design_dhw (column R) values for the first few building
# archetypes to the new capped figure, then leave the selection where the
# author left off (cell R6) as shown in the saved worksheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Value = 10000
$ws.Range("R4").Value = 10000
$ws.Range("R5").Value = 10000

$ws.Range("R6").Select() | Out-Null
